$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4206.5
$ws.Range("I74").Value = 3797.5334
$ws.Range("J74").Value = 5082.857
$ws.Range("K74").Value = 3797.5334
$ws.Range("L74").Value = 5082.857
$ws.Range("M74").Value = -2861.5334
$ws.Range("N74").Value = -6954.857

# Row 77
$ws.Range("H77").Value = 4206.5
$ws.Range("I77").Value = 3797.5334
$ws.Range("J77").Value = 5082.857
$ws.Range("K77").Value = 18987.667
$ws.Range("L77").Value = 25414.285
$ws.Range("M77").Value = -14307.667
$ws.Range("N77").Value = -34774.285

# Row 80
$ws.Range("H80").Value = 8606.791999999999
$ws.Range("I80").Value = 5343.3
$ws.Range("J80").Value = 10937.857
$ws.Range("K80").Value = 16029.9
$ws.Range("L80").Value = 32813.571
$ws.Range("M80").Value = -15031.9
$ws.Range("N80").Value = -34809.571

# Row 83
$ws.Range("H83").Value = 8606.791999999999
$ws.Range("I83").Value = 5343.3
$ws.Range("J83").Value = 10937.857
$ws.Range("K83").Value = 48089.7
$ws.Range("L83").Value = 98440.713
$ws.Range("M83").Value = -43097.7
$ws.Range("N83").Value = -108424.713

# Row 120
$ws.Range("H120").Value = 38760
$ws.Range("J120").Value = 38760
$ws.Range("L120").Value = 38760
$ws.Range("N120").Value = -48436

# Row 136
$ws.Range("H136").Value = 28575.6
$ws.Range("J136").Value = 28575.6
$ws.Range("L136").Value = 28575.6
$ws.Range("N136").Value = -38775.6

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5804.254
$ws.Range("I32").Value = 3671.1755
$ws.Range("J32").Value = 17962.8
$ws.Range("K32").Value = 3671.1755
$ws.Range("L32").Value = 17962.8
$ws.Range("M32").Value = -3384.1755
$ws.Range("N32").Value = -18536.8

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

# Row 110
$ws.Range("H110").Value = 2980.9285
$ws.Range("I110").Value = 2980.9285
$ws.Range("K110").Value = 2980.9285
$ws.Range("M110").Value = -935.9285

# Row 132
$ws.Range("H132").Value = 2362.486
$ws.Range("I132").Value = 2075.568
$ws.Range("J132").Value = 2813.3572
$ws.Range("K132").Value = 6226.704000000001
$ws.Range("L132").Value = 8440.071599999999
$ws.Range("M132").Value = -3696.704000000001
$ws.Range("N132").Value = -13500.0716

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 78679.46000000001
$ws.Range("I99").Value = 112532.11
$ws.Range("J99").Value = 2511
$ws.Range("K99").Value = 112532.11
$ws.Range("L99").Value = 2511
$ws.Range("M99").Value = -111034.11
$ws.Range("N99").Value = -5507

# Row 126
$ws.Range("H126").Value = 78679.46000000001
$ws.Range("I126").Value = 112532.11
$ws.Range("J126").Value = 2511
$ws.Range("K126").Value = 337596.33
$ws.Range("L126").Value = 7533
$ws.Range("M126").Value = -335126.33
$ws.Range("N126").Value = -12473

# Row 132
$ws.Range("H132").Value = 2738.913
$ws.Range("I132").Value = 1132.3
$ws.Range("J132").Value = 3974.7693
$ws.Range("K132").Value = 3396.9
$ws.Range("L132").Value = 11924.3079
$ws.Range("M132").Value = -866.8999999999996
$ws.Range("N132").Value = -16984.3079

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 706.2857
$ws.Range("I17").Value = 512.5
$ws.Range("J17").Value = 964.6667
$ws.Range("K17").Value = 1537.5
$ws.Range("L17").Value = 2894.0001
$ws.Range("M17").Value = -1368.5
$ws.Range("N17").Value = -3232.0001

# Row 131
$ws.Range("H131").Value = 1368.6097
$ws.Range("J131").Value = 1413.1538
$ws.Range("L131").Value = 4239.4614
$ws.Range("N131").Value = -14319.4614

# Row 137
$ws.Range("H137").Value = 12790.7
$ws.Range("I137").Value = 1303.6364
$ws.Range("J137").Value = 26830.445
$ws.Range("K137").Value = 3910.9092
$ws.Range("L137").Value = 80491.33499999999
$ws.Range("M137").Value = 1189.0908
$ws.Range("N137").Value = -90691.33499999999

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 10000
$ws.Range("J35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("N35").Value = -10596

# Row 70
$ws.Range("H70").Value = 4453.2856
$ws.Range("I70").Value = 4193.35
$ws.Range("J70").Value = 4799.8667
$ws.Range("K70").Value = 4193.35
$ws.Range("L70").Value = 4799.8667
$ws.Range("M70").Value = -3923.35
$ws.Range("N70").Value = -5339.8667

# Row 73
$ws.Range("H73").Value = 4453.2856
$ws.Range("I73").Value = 4193.35
$ws.Range("J73").Value = 4799.8667
$ws.Range("K73").Value = 4193.35
$ws.Range("L73").Value = 4799.8667
$ws.Range("M73").Value = -3257.35
$ws.Range("N73").Value = -6671.8667

# Row 80
$ws.Range("H80").Value = 4295.477
$ws.Range("I80").Value = 4743.788
$ws.Range("J80").Value = 2950.5454
$ws.Range("K80").Value = 4743.788
$ws.Range("L80").Value = 2950.5454
$ws.Range("M80").Value = -3745.788
$ws.Range("N80").Value = -4946.5454

# Row 83
$ws.Range("H83").Value = 4295.477
$ws.Range("I83").Value = 4743.788
$ws.Range("J83").Value = 2950.5454
$ws.Range("K83").Value = 23718.94
$ws.Range("L83").Value = 14752.727
$ws.Range("M83").Value = -18726.94
$ws.Range("N83").Value = -24736.727

# Row 113
$ws.Range("H113").Value = 2149.625
$ws.Range("I113").Value = 1821.75
$ws.Range("K113").Value = 1821.75
$ws.Range("M113").Value = 348.25

# Row 132
$ws.Range("H132").Value = 4262.9067
$ws.Range("I132").Value = 4482.2104
$ws.Range("J132").Value = 4089.2917
$ws.Range("K132").Value = 13446.6312
$ws.Range("L132").Value = 12267.8751
$ws.Range("M132").Value = -10916.6312
$ws.Range("N132").Value = -17327.8751

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1581.4375
$ws.Range("I82").Value = 828
$ws.Range("J82").Value = 1923.909
$ws.Range("K82").Value = 828
$ws.Range("L82").Value = 1923.909
$ws.Range("M82").Value = -467
$ws.Range("N82").Value = -2645.909

# Row 85
$ws.Range("H85").Value = 1581.4375
$ws.Range("I85").Value = 828
$ws.Range("J85").Value = 1923.909
$ws.Range("K85").Value = 828
$ws.Range("L85").Value = 1923.909
$ws.Range("M85").Value = 420
$ws.Range("N85").Value = -4419.909

# Row 132
$ws.Range("H132").Value = 9776.333000000001
$ws.Range("I132").Value = 3606.0625
$ws.Range("K132").Value = 10818.1875
$ws.Range("M132").Value = -8288.1875

# Row 133
$ws.Range("H133").Value = 34417.5
$ws.Range("J133").Value = 34417.5
$ws.Range("L133").Value = 34417.5
$ws.Range("N133").Value = -39477.5

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 6572.4443
$ws.Range("I34").Value = 6076
$ws.Range("J34").Value = 6714.2856
$ws.Range("K34").Value = 6076
$ws.Range("L34").Value = 6714.2856
$ws.Range("M34").Value = -5873
$ws.Range("N34").Value = -7120.2856

# Row 132
$ws.Range("H132").Value = 2765.6562
$ws.Range("I132").Value = 1861.2778
$ws.Range("J132").Value = 3928.4285
$ws.Range("K132").Value = 5583.8334
$ws.Range("L132").Value = 11785.2855
$ws.Range("M132").Value = -3053.8334
$ws.Range("N132").Value = -16845.2855
